# Natmi following Dr Hou advice
# Rebuild the LR-pair result table to include "ECs" as an additional sending
# cluster (in addition to the existing "sCs" sending cluster), each paired
# with all three target clusters (ECs, FAPs, sCs), for the Npnt -> Itgb1
# ligand-receptor pair.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns: Sending cluster, Ligand symbol, Receptor symbol, Target cluster,
#          E..T = the various expression / specificity metrics
$data = @(
    @("ECs","Npnt","Itgb1","ECs", 2, 0.6666666666666666, 0.391327, 1.173981, 0.2020349760212956, 0.2020349760212956, 3, 1, 117.044563, 351.133689, 0.3245365645427815, 0.3245365645427815, 45.802697705101, 412.224279345909, 0.0655677370354345, 0.06556773703543449),
    @("ECs","Npnt","Itgb1","FAPs", 2, 0.6666666666666666, 0.391327, 1.173981, 0.2020349760212956, 0.2020349760212956, 3, 1, 101.5800373333333, 304.740112, 0.281657135515876, 0.281657135515876, 39.75101126954133, 357.759101425872, 0.05690459262017681, 0.05690459262017679),
    @("ECs","Npnt","Itgb1","sCs", 2, 0.6666666666666666, 0.391327, 1.173981, 0.2020349760212956, 0.2020349760212956, 3, 1, 142.0267893333333, 426.080368, 0.3938062999413425, 0.3938062999413425, 55.57891738944533, 500.210256505008, 0.07956264636568428, 0.07956264636568426),
    @("sCs","Npnt","Itgb1","ECs", 3, 1, 1.5456, 4.6368, 0.7979650239787044, 0.7979650239787044, 3, 1, 117.044563, 351.133689, 0.3245365645427815, 0.3245365645427815, 180.9040765728, 1628.1366891552, 0.258968827507347, 0.258968827507347),
    @("sCs","Npnt","Itgb1","FAPs", 3, 1, 1.5456, 4.6368, 0.7979650239787044, 0.7979650239787044, 3, 1, 101.5800373333333, 304.740112, 0.281657135515876, 0.281657135515876, 157.0021057024, 1413.0189513216, 0.2247525428956992, 0.2247525428956992),
    @("sCs","Npnt","Itgb1","sCs", 3, 1, 1.5456, 4.6368, 0.7979650239787044, 0.7979650239787044, 3, 1, 142.0267893333333, 426.080368, 0.3938062999413425, 0.3938062999413425, 219.5166055936, 1975.6494503424, 0.3142436535756583, 0.3142436535756583)
)

$startRow = 2
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $values = $data[$i]
    for ($c = 0; $c -lt $values.Length; $c++) {
        $ws.Cells.Item($row, $c + 1).Value = $values[$c]
    }
}
